$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colA = New-Object "object[,]" 17,1
$colA[0,0] = 'Australian A-League Men'
$colA[1,0] = 'Turkish 2 Lig'
$colA[2,0] = 'Turkish 2 Lig'
$colA[3,0] = 'Turkish 1 Lig'
$colA[4,0] = 'Turkish 1 Lig'
$colA[5,0] = 'Italian Serie A'
$colA[6,0] = 'Turkish 1 Lig'
$colA[7,0] = 'Italian Serie A'
$colA[8,0] = 'English Premier League'
$colA[9,0] = 'Portuguese Primeira Liga'
$colA[10,0] = 'Portuguese Primeira Liga'
$colA[11,0] = 'Turkish 1 Lig'
$colA[12,0] = 'English Premier League'
$colA[13,0] = 'Italian Serie A'
$colA[14,0] = 'Portuguese Primeira Liga'
$colA[15,0] = 'Italian Serie A'
$colA[16,0] = 'Portuguese Primeira Liga'
$ws.Range("A2:A18").Value = $colA

$ws.Cells.Item(17, 2).Value = "'2025-12-28"
$ws.Cells.Item(18, 2).Value = "'2025-12-28"

$data = New-Object "object[,]" 17,39
$data[0,0] = '05:00:00'
$data[0,1] = 'Melbourne City'
$data[0,2] = 'Perth Glory'
$data[0,3] = 1.66
$data[0,4] = 1.7
$data[0,5] = 5.6
$data[0,6] = 6.4
$data[0,7] = 4.1
$data[0,8] = 4.3
$data[0,9] = 1.38
$data[0,10] = 1.06
$data[0,11] = 4
$data[0,12] = 1.3
$data[0,13] = 2.04
$data[0,14] = 1.87
$data[0,15] = 1.4
$data[0,16] = 3.2
$data[0,17] = 1.87
$data[0,18] = 2
$data[0,19] = 1.18
$data[0,20] = 2.42
$data[0,21] = 16
$data[0,22] = 21
$data[0,23] = 48
$data[0,24] = 190
$data[0,25] = 8.8
$data[0,26] = 9.4
$data[0,27] = 23
$data[0,28] = 85
$data[0,29] = 10
$data[0,30] = 9.8
$data[0,31] = 21
$data[0,32] = 85
$data[0,33] = 16.5
$data[0,34] = 17.5
$data[0,35] = 36
$data[0,36] = 130
$data[0,37] = 9.8
$data[0,38] = 100
$data[1,0] = '07:00:00'
$data[1,1] = 'Kepez Belediyespor'
$data[1,2] = 'Batman Petrolspor'
$data[1,3] = 1.01
$data[1,4] = 1000
$data[1,5] = 1.01
$data[1,6] = 1000
$data[1,7] = 1.02
$data[1,8] = 950
$data[1,9] = 1.01
$data[1,10] = 1.01
$data[1,11] = 1.24
$data[1,12] = 1.21
$data[1,13] = 1.24
$data[1,14] = 1.21
$data[1,15] = 1.18
$data[1,16] = 1.21
$data[1,17] = 1.01
$data[1,18] = 1.01
$data[1,19] = 1.01
$data[1,20] = 1.01
$data[1,21] = 1000
$data[1,22] = 1000
$data[1,23] = 1000
$data[1,24] = 1000
$data[1,25] = 1000
$data[1,26] = 1000
$data[1,27] = 1000
$data[1,28] = 1000
$data[1,29] = 1000
$data[1,30] = 1000
$data[1,31] = 1000
$data[1,32] = 1000
$data[1,33] = 1000
$data[1,34] = 1000
$data[1,35] = 1000
$data[1,36] = 1000
$data[1,37] = 1000
$data[1,38] = 1000
$data[2,0] = '07:00:00'
$data[2,1] = 'Iskenderunspor'
$data[2,2] = 'Ankaragucu'
$data[2,3] = 1.01
$data[2,4] = 1000
$data[2,5] = 1.01
$data[2,6] = 1000
$data[2,7] = 1.02
$data[2,8] = 950
$data[2,9] = 1.01
$data[2,10] = 1.01
$data[2,11] = 1.15
$data[2,12] = 1.01
$data[2,13] = 1.15
$data[2,14] = 1.01
$data[2,15] = 1.07
$data[2,16] = 1.01
$data[2,17] = 1.01
$data[2,18] = 1.01
$data[2,19] = 1.01
$data[2,20] = 1.01
$data[2,21] = 1000
$data[2,22] = 1000
$data[2,23] = 1000
$data[2,24] = 1000
$data[2,25] = 1000
$data[2,26] = 1000
$data[2,27] = 1000
$data[2,28] = 1000
$data[2,29] = 1000
$data[2,30] = 1000
$data[2,31] = 1000
$data[2,32] = 1000
$data[2,33] = 1000
$data[2,34] = 1000
$data[2,35] = 1000
$data[2,36] = 1000
$data[2,37] = 1000
$data[2,38] = 1000
$data[3,0] = '07:30:00'
$data[3,1] = 'Serik Belediyespor'
$data[3,2] = 'Boluspor'
$data[3,3] = 6
$data[3,4] = 8.4
$data[3,5] = 1.46
$data[3,6] = 1.56
$data[3,7] = 4.6
$data[3,8] = 5.3
$data[3,9] = 0
$data[3,10] = 0
$data[3,11] = 0
$data[3,12] = 0
$data[3,13] = 2.16
$data[3,14] = 1.68
$data[3,15] = 0
$data[3,16] = 0
$data[3,17] = 0
$data[3,18] = 0
$data[3,19] = 0
$data[3,20] = 0
$data[3,21] = 0
$data[3,22] = 0
$data[3,23] = 0
$data[3,24] = 0
$data[3,25] = 0
$data[3,26] = 0
$data[3,27] = 0
$data[3,28] = 0
$data[3,29] = 0
$data[3,30] = 0
$data[3,31] = 0
$data[3,32] = 0
$data[3,33] = 0
$data[3,34] = 0
$data[3,35] = 0
$data[3,36] = 0
$data[3,37] = 0
$data[3,38] = 0
$data[4,0] = '07:30:00'
$data[4,1] = 'Keciorengucu'
$data[4,2] = 'Umraniyespor'
$data[4,3] = 1.29
$data[4,4] = 1.46
$data[4,5] = 3.1
$data[4,6] = 1000
$data[4,7] = 4.8
$data[4,8] = 950
$data[4,9] = 0
$data[4,10] = 0
$data[4,11] = 0
$data[4,12] = 0
$data[4,13] = 2.02
$data[4,14] = 1.57
$data[4,15] = 0
$data[4,16] = 0
$data[4,17] = 0
$data[4,18] = 0
$data[4,19] = 0
$data[4,20] = 0
$data[4,21] = 0
$data[4,22] = 0
$data[4,23] = 0
$data[4,24] = 0
$data[4,25] = 0
$data[4,26] = 0
$data[4,27] = 0
$data[4,28] = 0
$data[4,29] = 0
$data[4,30] = 0
$data[4,31] = 0
$data[4,32] = 0
$data[4,33] = 0
$data[4,34] = 0
$data[4,35] = 0
$data[4,36] = 0
$data[4,37] = 0
$data[4,38] = 0
$data[5,0] = '08:30:00'
$data[5,1] = 'AC Milan'
$data[5,2] = 'Verona'
$data[5,3] = 1.4
$data[5,4] = 1.41
$data[5,5] = 10
$data[5,6] = 10.5
$data[5,7] = 5.3
$data[5,8] = 5.4
$data[5,9] = 0
$data[5,10] = 1.06
$data[5,11] = 4.1
$data[5,12] = 1.3
$data[5,13] = 2.06
$data[5,14] = 1.91
$data[5,15] = 1.41
$data[5,16] = 3.25
$data[5,17] = 2.24
$data[5,18] = 1.77
$data[5,19] = 0
$data[5,20] = 0
$data[5,21] = 17.5
$data[5,22] = 29
$data[5,23] = 100
$data[5,24] = 430
$data[5,25] = 7.6
$data[5,26] = 11.5
$data[5,27] = 38
$data[5,28] = 240
$data[5,29] = 7.6
$data[5,30] = 10.5
$data[5,31] = 34
$data[5,32] = 1000
$data[5,33] = 12
$data[5,34] = 15.5
$data[5,35] = 44
$data[5,36] = 240
$data[5,37] = 7
$data[5,38] = 1000
$data[6,0] = '10:00:00'
$data[6,1] = 'Erzurum BB'
$data[6,2] = 'Corum Belediyespor'
$data[6,3] = 2.12
$data[6,4] = 2.82
$data[6,5] = 2.88
$data[6,6] = 4.2
$data[6,7] = 3.05
$data[6,8] = 5.6
$data[6,9] = 0
$data[6,10] = 0
$data[6,11] = 0
$data[6,12] = 0
$data[6,13] = 1.64
$data[6,14] = 1.93
$data[6,15] = 0
$data[6,16] = 0
$data[6,17] = 0
$data[6,18] = 0
$data[6,19] = 0
$data[6,20] = 0
$data[6,21] = 0
$data[6,22] = 0
$data[6,23] = 0
$data[6,24] = 0
$data[6,25] = 0
$data[6,26] = 0
$data[6,27] = 0
$data[6,28] = 0
$data[6,29] = 0
$data[6,30] = 0
$data[6,31] = 0
$data[6,32] = 0
$data[6,33] = 0
$data[6,34] = 0
$data[6,35] = 0
$data[6,36] = 0
$data[6,37] = 0
$data[6,38] = 0
$data[7,0] = '11:00:00'
$data[7,1] = 'US Cremonese'
$data[7,2] = 'Napoli'
$data[7,3] = 7
$data[7,4] = 7.6
$data[7,5] = 1.6
$data[7,6] = 1.62
$data[7,7] = 4
$data[7,8] = 4.2
$data[7,9] = 0
$data[7,10] = 1.08
$data[7,11] = 3.4
$data[7,12] = 1.38
$data[7,13] = 1.81
$data[7,14] = 2.14
$data[7,15] = 1.3
$data[7,16] = 4
$data[7,17] = 2.18
$data[7,18] = 1.77
$data[7,19] = 0
$data[7,20] = 0
$data[7,21] = 13
$data[7,22] = 7
$data[7,23] = 8.6
$data[7,24] = 15
$data[7,25] = 20
$data[7,26] = 9.4
$data[7,27] = 11
$data[7,28] = 20
$data[7,29] = 65
$data[7,30] = 29
$data[7,31] = 29
$data[7,32] = 55
$data[7,33] = 1000
$data[7,34] = 160
$data[7,35] = 150
$data[7,36] = 1000
$data[7,37] = 1000
$data[7,38] = 11
$data[8,0] = '11:00:00'
$data[8,1] = 'Sunderland'
$data[8,2] = 'Leeds'
$data[8,3] = 2.74
$data[8,4] = 2.78
$data[8,5] = 2.92
$data[8,6] = 2.94
$data[8,7] = 3.35
$data[8,8] = 3.4
$data[8,9] = 0
$data[8,10] = 1.09
$data[8,11] = 3.2
$data[8,12] = 1.43
$data[8,13] = 1.75
$data[8,14] = 2.28
$data[8,15] = 1.27
$data[8,16] = 4.3
$data[8,17] = 1.96
$data[8,18] = 1.98
$data[8,19] = 0
$data[8,20] = 0
$data[8,21] = 11
$data[8,22] = 10
$data[8,23] = 19
$data[8,24] = 55
$data[8,25] = 9.4
$data[8,26] = 7.4
$data[8,27] = 13.5
$data[8,28] = 36
$data[8,29] = 16.5
$data[8,30] = 12.5
$data[8,31] = 20
$data[8,32] = 65
$data[8,33] = 42
$data[8,34] = 32
$data[8,35] = 55
$data[8,36] = 120
$data[8,37] = 34
$data[8,38] = 42
$data[9,0] = '12:30:00'
$data[9,1] = 'Arouca'
$data[9,2] = 'Gil Vicente'
$data[9,3] = 4.4
$data[9,4] = 4.7
$data[9,5] = 2.06
$data[9,6] = 2.1
$data[9,7] = 3.25
$data[9,8] = 3.45
$data[9,9] = 0
$data[9,10] = 0
$data[9,11] = 0
$data[9,12] = 0
$data[9,13] = 1.67
$data[9,14] = 2.28
$data[9,15] = 0
$data[9,16] = 0
$data[9,17] = 0
$data[9,18] = 0
$data[9,19] = 0
$data[9,20] = 0
$data[9,21] = 0
$data[9,22] = 0
$data[9,23] = 0
$data[9,24] = 0
$data[9,25] = 0
$data[9,26] = 0
$data[9,27] = 0
$data[9,28] = 0
$data[9,29] = 0
$data[9,30] = 0
$data[9,31] = 0
$data[9,32] = 0
$data[9,33] = 0
$data[9,34] = 0
$data[9,35] = 0
$data[9,36] = 0
$data[9,37] = 0
$data[9,38] = 0
$data[10,0] = '12:30:00'
$data[10,1] = 'Casa Pia'
$data[10,2] = 'Guimaraes'
$data[10,3] = 3.55
$data[10,4] = 3.85
$data[10,5] = 2.46
$data[10,6] = 2.58
$data[10,7] = 3
$data[10,8] = 3.1
$data[10,9] = 0
$data[10,10] = 0
$data[10,11] = 0
$data[10,12] = 0
$data[10,13] = 1.55
$data[10,14] = 2.48
$data[10,15] = 0
$data[10,16] = 0
$data[10,17] = 0
$data[10,18] = 0
$data[10,19] = 0
$data[10,20] = 0
$data[10,21] = 0
$data[10,22] = 0
$data[10,23] = 0
$data[10,24] = 0
$data[10,25] = 0
$data[10,26] = 0
$data[10,27] = 0
$data[10,28] = 0
$data[10,29] = 0
$data[10,30] = 0
$data[10,31] = 0
$data[10,32] = 0
$data[10,33] = 0
$data[10,34] = 0
$data[10,35] = 0
$data[10,36] = 0
$data[10,37] = 0
$data[10,38] = 0
$data[11,0] = '13:00:00'
$data[11,1] = 'Amed Sportif Faaliyetler'
$data[11,2] = '76 Igdir Belediyespor'
$data[11,3] = 1.82
$data[11,4] = 2.32
$data[11,5] = 3.35
$data[11,6] = 4.9
$data[11,7] = 3.4
$data[11,8] = 6.8
$data[11,9] = 0
$data[11,10] = 0
$data[11,11] = 0
$data[11,12] = 0
$data[11,13] = 1.96
$data[11,14] = 1.65
$data[11,15] = 0
$data[11,16] = 0
$data[11,17] = 0
$data[11,18] = 0
$data[11,19] = 0
$data[11,20] = 0
$data[11,21] = 0
$data[11,22] = 0
$data[11,23] = 0
$data[11,24] = 0
$data[11,25] = 0
$data[11,26] = 0
$data[11,27] = 0
$data[11,28] = 0
$data[11,29] = 0
$data[11,30] = 0
$data[11,31] = 0
$data[11,32] = 0
$data[11,33] = 0
$data[11,34] = 0
$data[11,35] = 0
$data[11,36] = 0
$data[11,37] = 0
$data[11,38] = 0
$data[12,0] = '13:30:00'
$data[12,1] = 'Crystal Palace'
$data[12,2] = 'Tottenham'
$data[12,3] = 2.34
$data[12,4] = 2.38
$data[12,5] = 3.4
$data[12,6] = 3.45
$data[12,7] = 3.45
$data[12,8] = 3.5
$data[12,9] = 0
$data[12,10] = 1.08
$data[12,11] = 3.75
$data[12,12] = 1.34
$data[12,13] = 1.9
$data[12,14] = 2.08
$data[12,15] = 1.35
$data[12,16] = 3.75
$data[12,17] = 1.83
$data[12,18] = 2.14
$data[12,19] = 0
$data[12,20] = 0
$data[12,21] = 13
$data[12,22] = 13
$data[12,23] = 24
$data[12,24] = 65
$data[12,25] = 10.5
$data[12,26] = 7.6
$data[12,27] = 15
$data[12,28] = 36
$data[12,29] = 15
$data[12,30] = 12
$data[12,31] = 21
$data[12,32] = 60
$data[12,33] = 34
$data[12,34] = 26
$data[12,35] = 42
$data[12,36] = 85
$data[12,37] = 20
$data[12,38] = 36
$data[13,0] = '14:00:00'
$data[13,1] = 'Bologna'
$data[13,2] = 'Sassuolo'
$data[13,3] = 1.76
$data[13,4] = 1.78
$data[13,5] = 5.7
$data[13,6] = 5.8
$data[13,7] = 3.85
$data[13,8] = 3.9
$data[13,9] = 0
$data[13,10] = 1.08
$data[13,11] = 3.65
$data[13,12] = 1.35
$data[13,13] = 1.9
$data[13,14] = 2.06
$data[13,15] = 1.34
$data[13,16] = 3.75
$data[13,17] = 1.98
$data[13,18] = 1.97
$data[13,19] = 0
$data[13,20] = 0
$data[13,21] = 14
$data[13,22] = 17.5
$data[13,23] = 44
$data[13,24] = 170
$data[13,25] = 8
$data[13,26] = 8.4
$data[13,27] = 22
$data[13,28] = 85
$data[13,29] = 10.5
$data[13,30] = 10
$data[13,31] = 22
$data[13,32] = 85
$data[13,33] = 18.5
$data[13,34] = 19.5
$data[13,35] = 42
$data[13,36] = 150
$data[13,37] = 12.5
$data[13,38] = 120
$data[14,0] = '15:00:00'
$data[14,1] = 'Braga'
$data[14,2] = 'Benfica'
$data[14,3] = 3.55
$data[14,4] = 3.95
$data[14,5] = 2.24
$data[14,6] = 2.36
$data[14,7] = 3.3
$data[14,8] = 3.55
$data[14,9] = 0
$data[14,10] = 0
$data[14,11] = 0
$data[14,12] = 0
$data[14,13] = 1.78
$data[14,14] = 2.12
$data[14,15] = 0
$data[14,16] = 0
$data[14,17] = 0
$data[14,18] = 0
$data[14,19] = 0
$data[14,20] = 0
$data[14,21] = 0
$data[14,22] = 0
$data[14,23] = 0
$data[14,24] = 0
$data[14,25] = 0
$data[14,26] = 0
$data[14,27] = 0
$data[14,28] = 0
$data[14,29] = 0
$data[14,30] = 0
$data[14,31] = 0
$data[14,32] = 0
$data[14,33] = 0
$data[14,34] = 0
$data[14,35] = 0
$data[14,36] = 0
$data[14,37] = 0
$data[14,38] = 0
$data[15,0] = '16:45:00'
$data[15,1] = 'Atalanta'
$data[15,2] = 'Inter'
$data[15,3] = 4
$data[15,4] = 4.1
$data[15,5] = 2.08
$data[15,6] = 2.1
$data[15,7] = 3.65
$data[15,8] = 3.75
$data[15,9] = 0
$data[15,10] = 1.06
$data[15,11] = 4.4
$data[15,12] = 1.27
$data[15,13] = 2.16
$data[15,14] = 1.82
$data[15,15] = 1.46
$data[15,16] = 3.05
$data[15,17] = 1.7
$data[15,18] = 2.32
$data[15,19] = 0
$data[15,20] = 0
$data[15,21] = 17.5
$data[15,22] = 11.5
$data[15,23] = 14
$data[15,24] = 26
$data[15,25] = 16.5
$data[15,26] = 8.4
$data[15,27] = 11
$data[15,28] = 21
$data[15,29] = 30
$data[15,30] = 16
$data[15,31] = 16.5
$data[15,32] = 34
$data[15,33] = 75
$data[15,34] = 44
$data[15,35] = 50
$data[15,36] = 80
$data[15,37] = 40
$data[15,38] = 13.5
$data[16,0] = '17:30:00'
$data[16,1] = 'Sporting Lisbon'
$data[16,2] = 'Rio Ave'
$data[16,3] = 1.16
$data[16,4] = 1.17
$data[16,5] = 21
$data[16,6] = 26
$data[16,7] = 9.4
$data[16,8] = 10
$data[16,9] = 0
$data[16,10] = 0
$data[16,11] = 0
$data[16,12] = 0
$data[16,13] = 3.1
$data[16,14] = 1.39
$data[16,15] = 0
$data[16,16] = 0
$data[16,17] = 0
$data[16,18] = 0
$data[16,19] = 0
$data[16,20] = 0
$data[16,21] = 0
$data[16,22] = 0
$data[16,23] = 0
$data[16,24] = 0
$data[16,25] = 0
$data[16,26] = 0
$data[16,27] = 0
$data[16,28] = 0
$data[16,29] = 0
$data[16,30] = 0
$data[16,31] = 0
$data[16,32] = 0
$data[16,33] = 0
$data[16,34] = 0
$data[16,35] = 0
$data[16,36] = 0
$data[16,37] = 0
$data[16,38] = 0
$ws.Range("C2:AO18").Value = $data
